# Update "想去人数" (column F) figures across all four sheets to match the
# freshly generated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 649
$ws.Cells.Item(4, 6).Value = 812
$ws.Cells.Item(5, 6).Value = 496
$ws.Cells.Item(6, 6).Value = 388
$ws.Cells.Item(7, 6).Value = 475
$ws.Cells.Item(8, 6).Value = 884
$ws.Cells.Item(10, 6).Value = 806
$ws.Cells.Item(11, 6).Value = 636
$ws.Cells.Item(12, 6).Value = 107
$ws.Cells.Item(15, 6).Value = 716
$ws.Cells.Item(16, 6).Value = 205
$ws.Cells.Item(17, 6).Value = 507
$ws.Cells.Item(18, 6).Value = 456
$ws.Cells.Item(19, 6).Value = 1221
$ws.Cells.Item(20, 6).Value = 103
$ws.Cells.Item(21, 6).Value = 920
$ws.Cells.Item(22, 6).Value = 2667
$ws.Cells.Item(23, 6).Value = 1146
$ws.Cells.Item(24, 6).Value = 610
$ws.Cells.Item(25, 6).Value = 137
$ws.Cells.Item(26, 6).Value = 1194
$ws.Cells.Item(27, 6).Value = 48
$ws.Cells.Item(28, 6).Value = 885
$ws.Cells.Item(29, 6).Value = 88
$ws.Cells.Item(30, 6).Value = 1221

# --- Sheet "演出" (Shows) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 475

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 693

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 693
$ws.Cells.Item(4, 6).Value = 649
$ws.Cells.Item(5, 6).Value = 812
$ws.Cells.Item(6, 6).Value = 496
$ws.Cells.Item(8, 6).Value = 388
$ws.Cells.Item(9, 6).Value = 475
$ws.Cells.Item(10, 6).Value = 475
$ws.Cells.Item(14, 6).Value = 884
$ws.Cells.Item(16, 6).Value = 806
$ws.Cells.Item(17, 6).Value = 636
$ws.Cells.Item(18, 6).Value = 107
$ws.Cells.Item(26, 6).Value = 716
$ws.Cells.Item(27, 6).Value = 205
$ws.Cells.Item(28, 6).Value = 507
$ws.Cells.Item(29, 6).Value = 456
$ws.Cells.Item(30, 6).Value = 1221
$ws.Cells.Item(31, 6).Value = 103
$ws.Cells.Item(32, 6).Value = 920
$ws.Cells.Item(33, 6).Value = 2667
$ws.Cells.Item(34, 6).Value = 1146
$ws.Cells.Item(35, 6).Value = 610
$ws.Cells.Item(36, 6).Value = 137
$ws.Cells.Item(37, 6).Value = 1194
$ws.Cells.Item(38, 6).Value = 48
$ws.Cells.Item(40, 6).Value = 885
$ws.Cells.Item(41, 6).Value = 88
$ws.Cells.Item(42, 6).Value = 1221
